# Regenerate save_data: replace column G ("K", formerly Strike#) values
# with newly computed K values (per commit message: "regen save_data to
# use K instead of Strike#, regen std/mean, calc and write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (header "K"), rows 2-58, in row order.
$newK = @(0,2,1,0,1,2,1,1,1,2,1,0,0,4,1,1,3,0,2,2,2,3,0,2,0,0,0,2,0,1,0,0,1,1,2,0,1,1,1,0,1,0,0,3,1,0,1,1,0,2,0,2,0,1,1,2,1)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
